$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = -3336
# Row 98
$ws.Range("H98").Value = 3899.8333
$ws.Range("I98").Value = 4209.25
$ws.Range("J98").Value = 1424.5
$ws.Range("K98").Value = 4209.25
$ws.Range("L98").Value = 1424.5
$ws.Range("M98").Value = -2711.25
$ws.Range("N98").Value = -4420.5
# Row 107
$ws.Range("H107").Value = 3049.1482
$ws.Range("I107").Value = 2382.9524
$ws.Range("J107").Value = 5380.8335
$ws.Range("K107").Value = 2382.9524
$ws.Range("L107").Value = 5380.8335
$ws.Range("M107").Value = -462.9524000000001
$ws.Range("N107").Value = -9220.833500000001
# Row 116
$ws.Range("H116").Value = 3188.9092
$ws.Range("I116").Value = 3230.889
$ws.Range("K116").Value = 3230.889
$ws.Range("M116").Value = 211.1109999999999
# Row 122
$ws.Range("H122").Value = 3899.8333
$ws.Range("I122").Value = 4209.25
$ws.Range("J122").Value = 1424.5
$ws.Range("K122").Value = 12627.75
$ws.Range("L122").Value = 4273.5
$ws.Range("M122").Value = -10177.75
$ws.Range("N122").Value = -9173.5
# Row 138
$ws.Range("H138").Value = 2660.7073
$ws.Range("J138").Value = 2599.7068
$ws.Range("L138").Value = 7799.1204
$ws.Range("N138").Value = -18079.1204

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1250
$ws.Range("I2").Value = 1250
$ws.Range("K2").Value = 1250
$ws.Range("M2").Value = -1137
# Row 32
$ws.Range("H32").Value = 10350.655
$ws.Range("I32").Value = 7735.086
$ws.Range("J32").Value = 19505.15
$ws.Range("K32").Value = 7735.086
$ws.Range("L32").Value = 19505.15
$ws.Range("M32").Value = -7448.086
$ws.Range("N32").Value = -20079.15
# Row 45
$ws.Range("H45").Value = 1324.5555
$ws.Range("I45").Value = 1120.1666
$ws.Range("K45").Value = 1120.1666
$ws.Range("M45").Value = -743.1666
# Row 61
$ws.Range("H61").Value = 142859310
$ws.Range("I61").Value = 250001550
$ws.Range("J61").Value = 2999.6667
$ws.Range("K61").Value = 250001550
$ws.Range("L61").Value = 2999.6667
$ws.Range("M61").Value = -250001338
$ws.Range("N61").Value = -3423.6667
# Row 74
$ws.Range("H74").Value = 1368.5518
$ws.Range("I74").Value = 963.6818
$ws.Range("J74").Value = 2641
$ws.Range("K74").Value = 963.6818
$ws.Range("L74").Value = 2641
$ws.Range("M74").Value = -89.68179999999995
$ws.Range("N74").Value = -4389
# Row 77
$ws.Range("H77").Value = 1368.5518
$ws.Range("I77").Value = 963.6818
$ws.Range("J77").Value = 2641
$ws.Range("K77").Value = 4818.409
$ws.Range("L77").Value = 13205
$ws.Range("M77").Value = -450.4089999999997
$ws.Range("N77").Value = -21941
# Row 110
$ws.Range("H110").Value = 280.81818
$ws.Range("I110").Value = 301.5
$ws.Range("K110").Value = 301.5
$ws.Range("M110").Value = 1743.5
# Row 116
$ws.Range("H116").Value = 1250
$ws.Range("I116").Value = 1250
$ws.Range("K116").Value = 1250
$ws.Range("M116").Value = 1044
# Row 125
$ws.Range("H125").Value = 36225
$ws.Range("J125").Value = 36225
$ws.Range("L125").Value = 36225
$ws.Range("N125").Value = -46065
# Row 132
$ws.Range("H132").Value = 3750.5173
$ws.Range("I132").Value = 3470.5908
$ws.Range("J132").Value = 4630.2856
$ws.Range("K132").Value = 10411.7724
$ws.Range("L132").Value = 13890.8568
$ws.Range("M132").Value = -7881.7724
$ws.Range("N132").Value = -18950.8568
# Row 136
$ws.Range("H136").Value = 142859310
$ws.Range("I136").Value = 250001550
$ws.Range("J136").Value = 2999.6667
$ws.Range("K136").Value = 750004650
$ws.Range("L136").Value = 8999.000100000001
$ws.Range("M136").Value = -750002100
$ws.Range("N136").Value = -14099.0001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1250
$ws.Range("I3").Value = 1250
$ws.Range("K3").Value = 1250
$ws.Range("M3").Value = -1136
# Row 22
$ws.Range("H22").Value = 475.1111
$ws.Range("I22").Value = 449.75
$ws.Range("K22").Value = 449.75
$ws.Range("M22").Value = -276.75
# Row 86
$ws.Range("H86").Value = 4193.8887
$ws.Range("I86").Value = 4291.077
$ws.Range("K86").Value = 4291.077
$ws.Range("M86").Value = -3168.077
# Row 89
$ws.Range("H89").Value = 4193.8887
$ws.Range("I89").Value = 4291.077
$ws.Range("K89").Value = 21455.385
$ws.Range("M89").Value = -15839.385
# Row 107
$ws.Range("H107").Value = 924.86664
$ws.Range("I107").Value = 846.6667
$ws.Range("K107").Value = 846.6667
$ws.Range("M107").Value = 1073.3333
# Row 134
$ws.Range("H134").Value = 4471.793
$ws.Range("I134").Value = 1035.5416
$ws.Range("K134").Value = 3106.6248
$ws.Range("M134").Value = -571.6248000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 111112824
$ws.Range("I16").Value = 111112824
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 111112824
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -111112537
$ws.Range("N16").Value = ""
# Row 31
$ws.Range("H31").Value = 1536.3846
$ws.Range("I31").Value = 1367.5
$ws.Range("J31").Value = 2465.25
$ws.Range("K31").Value = 1367.5
$ws.Range("L31").Value = 2465.25
$ws.Range("M31").Value = -1072.5
$ws.Range("N31").Value = -3055.25
# Row 34
$ws.Range("H34").Value = 1536.3846
$ws.Range("I34").Value = 1367.5
$ws.Range("J34").Value = 2465.25
$ws.Range("K34").Value = 1367.5
$ws.Range("L34").Value = 2465.25
$ws.Range("M34").Value = -1165.5
$ws.Range("N34").Value = -2869.25
# Row 86
$ws.Range("H86").Value = 3200840
$ws.Range("I86").Value = 5145972
$ws.Range("J86").Value = 40000.625
$ws.Range("K86").Value = 5145972
$ws.Range("L86").Value = 40000.625
$ws.Range("M86").Value = -5144849
$ws.Range("N86").Value = -42246.625
# Row 89
$ws.Range("H89").Value = 3200840
$ws.Range("I89").Value = 5145972
$ws.Range("J89").Value = 40000.625
$ws.Range("K89").Value = 25729860
$ws.Range("L89").Value = 200003.125
$ws.Range("M89").Value = -25724244
$ws.Range("N89").Value = -211235.125
# Row 107
$ws.Range("H107").Value = 679.24
$ws.Range("I107").Value = 339.8
$ws.Range("J107").Value = 2037
$ws.Range("K107").Value = 339.8
$ws.Range("L107").Value = 2037
$ws.Range("M107").Value = 1580.2
$ws.Range("N107").Value = -5877
# Row 113
$ws.Range("H113").Value = 111112824
$ws.Range("I113").Value = 111112824
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 111112824
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -111110654
$ws.Range("N113").Value = ""
# Row 122
$ws.Range("H122").Value = 1198.1904
$ws.Range("I122").Value = 815.6
$ws.Range("J122").Value = 2154.6667
$ws.Range("K122").Value = 2446.8
$ws.Range("L122").Value = 6464.000100000001
$ws.Range("M122").Value = 3.199999999999818
$ws.Range("N122").Value = -11364.0001
# Row 132
$ws.Range("H132").Value = 1750.3636
$ws.Range("I132").Value = 1514.6818
$ws.Range("J132").Value = 2221.7273
$ws.Range("K132").Value = 4544.0454
$ws.Range("L132").Value = 6665.1819
$ws.Range("M132").Value = -2014.0454
$ws.Range("N132").Value = -11725.1819

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 691.5
$ws.Range("I113").Value = 604.7222
$ws.Range("K113").Value = 1814.1666
$ws.Range("M113").Value = 355.8334
# Row 131
$ws.Range("H131").Value = 21309420
$ws.Range("J131").Value = 46591.91
$ws.Range("L131").Value = 139775.73
$ws.Range("N131").Value = -149855.73
# Row 132
$ws.Range("H132").Value = 1349.0834
$ws.Range("I132").Value = 1056
$ws.Range("J132").Value = 1558.4286
$ws.Range("K132").Value = 9504
$ws.Range("L132").Value = 14025.8574
$ws.Range("M132").Value = -6974
$ws.Range("N132").Value = -19085.8574

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 352.1154
$ws.Range("I2").Value = 239.6923
$ws.Range("J2").Value = 464.53845
$ws.Range("K2").Value = 239.6923
$ws.Range("L2").Value = 464.53845
$ws.Range("M2").Value = -126.6923
$ws.Range("N2").Value = -690.53845
# Row 113
$ws.Range("H113").Value = 1177.8572
$ws.Range("I113").Value = 1222.2222
$ws.Range("J113").Value = 1098
$ws.Range("K113").Value = 1222.2222
$ws.Range("L113").Value = 1098
$ws.Range("M113").Value = 947.7778000000001
$ws.Range("N113").Value = -5438
# Row 122
$ws.Range("H122").Value = 1199.8334
$ws.Range("I122").Value = 1099.75
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 3299.25
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -849.25
$ws.Range("N122").Value = -9100
# Row 132
$ws.Range("H132").Value = 5422.8857
$ws.Range("I132").Value = 6832.476
$ws.Range("J132").Value = 3308.5
$ws.Range("K132").Value = 20497.428
$ws.Range("L132").Value = 9925.5
$ws.Range("M132").Value = -17967.428
$ws.Range("N132").Value = -14985.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2013.8
$ws.Range("I7").Value = 1620.7
$ws.Range("J7").Value = 2406.9
$ws.Range("K7").Value = 1620.7
$ws.Range("L7").Value = 2406.9
$ws.Range("M7").Value = -1508.7
$ws.Range("N7").Value = -2630.9
# Row 61
$ws.Range("H61").Value = 1071.9565
$ws.Range("I61").Value = 917.2632
$ws.Range("K61").Value = 917.2632
$ws.Range("M61").Value = -715.2632
# Row 113
$ws.Range("H113").Value = 1071.9565
$ws.Range("I113").Value = 917.2632
$ws.Range("K113").Value = 917.2632
$ws.Range("M113").Value = 1252.7368
# Row 122
$ws.Range("H122").Value = 70834456
$ws.Range("I122").Value = 94444940
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 283334820
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -283332370
$ws.Range("N122").Value = -13900
# Row 126
$ws.Range("H126").Value = 2013.8
$ws.Range("I126").Value = 1620.7
$ws.Range("J126").Value = 2406.9
$ws.Range("K126").Value = 4862.1
$ws.Range("L126").Value = 7220.700000000001
$ws.Range("M126").Value = -2392.1
$ws.Range("N126").Value = -12160.7

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 529.8333
$ws.Range("I107").Value = 529.8333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1589.4999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 330.5001
$ws.Range("N107").Value = ""
# Row 132
$ws.Range("H132").Value = 1854.3276
$ws.Range("I132").Value = 1718.96
$ws.Range("K132").Value = 5156.88
$ws.Range("M132").Value = -2626.88
